$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5; existing rows 5-12 shift down to 6-13.
$ws.Rows(5).Insert()

# Match the date number-format used by the surrounding rows in column D.
$ws.Range("D5").NumberFormat = $ws.Range("D6").NumberFormat

# Populate the new row 5 with the new data record.
$ws.Range("A5").Value = 11
$ws.Range("B5").Value = "Vega Monumental Concepción"
$ws.Range("C5").Value = "Bíobío"
$ws.Range("D5").Value = 44482
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 100112026
$ws.Range("G5").Value = "Haba"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 430
$ws.Range("K5").Value = 8000
$ws.Range("L5").Value = 8500
$ws.Range("M5").Value = 8267
$ws.Range("N5").Value = "`$/saco 25 kilos"
$ws.Range("O5").Value = "Región de O'Higgins"
$ws.Range("P5").Value = 331
$ws.Range("Q5").Value = 25
$ws.Range("R5").Value = "Hortaliza"
